# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) values in the 展览 (Exhibition)
# and 全部类型 (All types) sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8833
$ws1.Range("F16").Value = 8727
$ws1.Range("F26").Value = 67
$ws1.Range("F34").Value = 2202
$ws1.Range("F35").Value = 868
$ws1.Range("F40").Value = 242
$ws1.Range("F41").Value = 172
$ws1.Range("F43").Value = 567
$ws1.Range("F44").Value = 79

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8833
$ws4.Range("F20").Value = 8727
$ws4.Range("F28").Value = 67
$ws4.Range("F34").Value = 2202
$ws4.Range("F35").Value = 868
$ws4.Range("F40").Value = 242
$ws4.Range("F42").Value = 172
